$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Corrected Avg_Agent_Step_Time (G), Avg_Experiment_Time (H),
# Std_Agent_Step_Time (M), Std_Experiment_Time (N) values.
$values = @{
    2  = @{ G = 2.46580998;         H = 169.88336374;        M = 0.2575856796067603;  N = 35.65025453032664 }
    3  = @{ G = 1.51946332;         H = 171.45638814;        M = 0.4563816235617154;  N = 60.62167444470195 }
    4  = @{ G = 1.31476032;         H = 60.33712061999999;   M = 0.2901926646659931;  N = 22.75476803659272 }
    5  = @{ G = 0.7536861799999999; H = 48.68516868;         M = 0.2075764034850184;  N = 17.96604165128577 }
    6  = @{ G = 0.5826141199999999; H = 16.24644572;         M = 0.1705166284447882;  N = 7.710043693630528 }
    7  = @{ G = 0.3537025600000001; H = 12.07943624;         M = 0.1064298320441421;  N = 4.768850306221507 }
    8  = @{ G = 0.42461368;         H = 4.05563942;          M = 0.1572723374045257;  N = 2.259023269990147 }
    9  = @{ G = 0.47369198;         H = 8.272052520000001;   M = 0.1169436231843642;  N = 3.776438336769833 }
    10 = @{ G = 0.2038068;          H = 1.46559084;          M = 0.08619550150860913; N = 0.9145212355313407 }
    11 = @{ G = 0.26120484;         H = 3.65293844;          M = 0.08015625125843905; N = 2.300968590328717 }
    12 = @{ G = 0.13002958;         H = 0.80684384;          M = 0.06834593545933235; N = 0.6244658429654223 }
    13 = @{ G = 0.14244442;         H = 1.61395626;          M = 0.04614002941581838; N = 0.940494734128831 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("G$row").Value = $rowVals.G
    $ws.Range("H$row").Value = $rowVals.H
    $ws.Range("M$row").Value = $rowVals.M
    $ws.Range("N$row").Value = $rowVals.N
}
